$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows, per the diff.
$updates = @{
    3  = -3
    4  = -4
    5  = 1
    7  = 4
    8  = 2
    9  = -2
    10 = 12
    12 = 2
    13 = 5
    14 = 5
    15 = -1
    16 = 5
    17 = -2
    18 = -4
    19 = 7
    20 = 0
    21 = -2
    22 = -2
    24 = 2
    25 = -2
    26 = 4
    27 = 9
    28 = 7
    29 = -2
    30 = -2
    31 = 4
    32 = -4
    33 = 4
    34 = -2
    35 = -2
    36 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
